$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 75
$ws.Range("F6").Value = 2725
$ws.Range("F9").Value = 7393
$ws.Range("F11").Value = 7571
$ws.Range("F14").Value = 4
$ws.Range("F15").Value = 6045
$ws.Range("F16").Value = 3232
$ws.Range("F17").Value = 3600
$ws.Range("F19").Value = 2
$ws.Range("F20").Value = 11
$ws.Range("F22").Value = 434
$ws.Range("F23").Value = 2
$ws.Range("F25").Value = 272
$ws.Range("F26").Value = 2085
$ws.Range("F30").Value = 253
$ws.Range("F31").Value = 1047
$ws.Range("F34").Value = 2584
$ws.Range("F35").Value = 1424
$ws.Range("F37").Value = 5
$ws.Range("F38").Value = 9
$ws.Range("F39").Value = 3181
$ws.Range("F40").Value = 142
$ws.Range("F41").Value = 234
$ws.Range("F44").Value = 472
$ws.Range("F45").Value = 1240

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 53
$ws.Range("F9").Value = 394
$ws.Range("F13").Value = 13
$ws.Range("F14").Value = 31

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 117

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 75
$ws.Range("F7").Value = 53
$ws.Range("F9").Value = 117
$ws.Range("F10").Value = 2725
$ws.Range("F14").Value = 7393
$ws.Range("F16").Value = 7571
$ws.Range("F18").Value = 6045
$ws.Range("F19").Value = 3232
$ws.Range("F20").Value = 3600
$ws.Range("F22").Value = 11
$ws.Range("F24").Value = 434
$ws.Range("F29").Value = 272
$ws.Range("F30").Value = 2085
$ws.Range("F31").Value = 13
$ws.Range("F33").Value = 31
$ws.Range("F36").Value = 1047
$ws.Range("F38").Value = 2584
$ws.Range("F39").Value = 1424
$ws.Range("F41").Value = 3181
$ws.Range("F42").Value = 142
$ws.Range("F45").Value = 472
$ws.Range("F46").Value = 1240
